# Update the "Marking"/"Total" rows of the marksheet with correct
# Right-answer counts and the corrected "correct/total marks" summary.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> Right count 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right count 60 -> 100
$ws.Range("B12").Value = 100

# Row 12 "Total" -> Max column (correct/total marks) "55/84" -> "100/140"
$ws.Range("E12").Value = "100/140"
